# Updates the SizeDecisionsSpreadsheet:
#  - inserts a new top banner row with a merged "AVERAGE DECISIONS MADE" cell
#  - removes the BCP (column D) data values from the size table
#  - moves the selection to F13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new banner row above the table -------------------
# Insert a row at the very top: this shifts the header (was row 1) to row 2
# and the size/data rows (were 2-11) down to 3-12. It also pushes the blank
# trailer rows (15-26) down to 16-27, so immediately remove the row that
# Insert() vacated right below the table (the empty row 13) to slide the
# untouched trailer rows (16-27) back up to their original numbers (15-26).
$ws.Rows("1:1").Insert()
$ws.Rows("13:13").Delete()

# --- 2. Populate the new banner row -----------------------------------------
$ws.Range("B1:F1").Merge()
$ws.Range("B1").Value = "AVERAGE DECISIONS MADE"
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").Borders.Item(9).LineStyle = -4119
$ws.Range("B1").Borders.Item(9).Color = 4144959
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)
$ws.Rows("1:1").RowHeight = 15.75

# --- 3. Drop the BCP (column D) values from the size table ------------------
$ws.Range("D3:D12").Clear()

# --- 4. Restore the active selection ----------------------------------------
$ws.Range("F13").Select()
